# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method

$wb = $excel.ActiveWorkbook

# Update "OFF" sheet, row 3 (Receiving/Rushing totals row labeled "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 229
$wsOff.Range("C3").Value = 164
$wsOff.Range("D3").Value = 56
$wsOff.Range("E3").Value = 28
$wsOff.Range("F3").Value = 3
$wsOff.Range("G3").Value = 4

# Update "DEF" sheet, row 3
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 205
$wsDef.Range("C3").Value = 150
$wsDef.Range("D3").Value = 39
$wsDef.Range("E3").Value = 16
